$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.414.88'
$ws.Range('D3').Value = '1.842.56'
$ws.Range('E3').Value = '  +1.96%  '
$ws.Range('D5').Value = '''229.60'
$ws.Range('D6').Value = '''0.610'
$ws.Range('E6').Value = '  +2.46%  '
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('D8').Value = '''43.33'
$ws.Range('E8').Value = '  +13.08%  '
$ws.Range('E9').Value = '  +7.15%  '
$ws.Range('E10').Value = '  +3.48%  '
$ws.Range('D12').Value = '2.109.25'
$ws.Range('E12').Value = '  +2.00%  '
$ws.Range('D13').Value = '1.838.48'
$ws.Range('E13').Value = '  +1.83%  '
$ws.Range('D14').Value = '''11.31'
$ws.Range('E14').Value = '  +2.06%  '
$ws.Range('D15').Value = '''0.674'
$ws.Range('E15').Value = '  +7.24%  '
$ws.Range('E16').Value = '  +7.02%  '
$ws.Range('D17').Value = '35.419.05'
$ws.Range('E17').Value = '  +2.83%  '
$ws.Range('E18').Value = '  +3.16%  '
$ws.Range('D19').Value = '0.0₃0795'
$ws.Range('E19').Value = '  +3.40%  '
$ws.Range('D20').Value = '''244.56'
$ws.Range('E20').Value = '  +1.19%  '
$ws.Range('D21').Value = '''12.11'
$ws.Range('E21').Value = '  +9.26%  '
$ws.Range('D22').Value = '''4.67'
$ws.Range('E22').Value = '  +13.87%  '
$ws.Range('E23').Value = '  +0.21%  '
$ws.Range('E24').Value = '  +0.83%  '
$ws.Range('D25').Value = '''168.93'
$ws.Range('E25').Value = '  -1.15%  '
$ws.Range('D26').Value = '''7.91'
$ws.Range('E26').Value = '  +2.50%  '
$ws.Range('D27').Value = '''17.78'
$ws.Range('E27').Value = '  +2.29%  '
$ws.Range('E28').Value = '  +1.84%  '
$ws.Range('E29').Value = '  +13.19%  '
$ws.Range('E30').Value = '  +0.26%  '
$ws.Range('D31').Value = '3.398.35'
$ws.Range('E31').Value = '  +39.87%  '
$ws.Range('E32').Value = '  +6.30%  '
$ws.Range('E34').Value = '  +4.30%  '
$ws.Range('E35').Value = '  +2.67%  '
$ws.Range('D36').Value = '''95.81'
$ws.Range('E36').Value = '  +15.56%  '
$ws.Range('D37').Value = '''0.690'
$ws.Range('E37').Value = '  +7.63%  '
$ws.Range('D38').Value = '1.346.24'
$ws.Range('E38').Value = '  +2.19%  '
$ws.Range('D39').Value = '''1.08'
$ws.Range('E39').Value = '  +2.87%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = '''0.0194'
$ws.Range('E40').Value = '  +3.59%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').Value = '''2.43'
$ws.Range('E41').Value = '  +5.51%  '
$ws.Range('E42').Value = '  +6.35%  '
$ws.Range('B43').Value = 'InjectiveProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D43').Value = '''14.99'
$ws.Range('E43').Value = '  +9.50%  '
$ws.Range('B44').Value = 'WEMIXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').Value = '''1.26'
$ws.Range('E44').Value = '  +3.75%  '
$ws.Range('E45').Value = '  +0.73%  '
$ws.Range('D46').Value = '''2.81'
$ws.Range('E46').Value = '  -0.19%  '
$ws.Range('D47').Value = '''6.23'
$ws.Range('E47').Value = '  +7.86%  '
$ws.Range('D48').Value = '''0.0520'
$ws.Range('E48').Value = '  +1.55%  '
$ws.Range('D49').Value = '2.009.89'
$ws.Range('E49').Value = '  +2.11%  '
$ws.Range('E50').Value = '  +0.25%  '
$ws.Range('D51').Value = '''102.90'
$ws.Range('E51').Value = '  +0.66%  '
